$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B unit-name values for the new rows (rows 81..124 -> A values 80..123)
$units = @(
  "Apothecary",
  "Apothecary",
  "Apothecary",
  "Heavy Intercessor",
  "Heavy Intercessor",
  "Heavy Intercessor",
  "Heavy Intercessor",
  "Terminator",
  "Captain In Terminator Armor",
  "Terminator",
  "Terminator",
  "Ballistus Dreadnought",
  "Ballistus Dreadnought",
  "Ballistus Dreadnought",
  "Ballistus Dreadnought",
  "Sternguard Veteran",
  "Sternguard Veteran",
  "Sternguard Veteran",
  "Sternguard Veteran",
  "Sternguard Veteran",
  "Sternguard Veteran",
  "Sternguard Veteran",
  "Hellblaster",
  "Hellblaster",
  "Hellblaster",
  "Techmarine",
  "Techmarine",
  "Assault Intercessor With Jump Pack",
  "Assault Intercessor With Jump Pack",
  "Assault Intercessor With Jump Pack",
  "Assault Intercessor With Jump Pack",
  "Infernus Marine",
  "Infernus Marine",
  "Agressor",
  "Agressor",
  "Agressor",
  "Ballistus Dreadnought",
  "Ballistus Dreadnought",
  "Ballistus Dreadnought",
  "Apothecary",
  "Apothecary",
  "Apothecary",
  "Apothecary",
  "Apothecary"
)

$startRow = 81
for ($i = 0; $i -lt $units.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 1).Value = $row - 1
  $ws.Cells.Item($row, 2).Value = $units[$i]
}

$ws.Range("E127").Select()
$excel.ActiveWindow.ScrollRow = 106
